$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.079490228806947
$ws.Range("C2").Value = 0.2879520647798586
$ws.Range("D2").Value = 0.06996198002806508
$ws.Range("E2").Value = 0.04325744321653202
$ws.Range("F2").Value = 1.816407904982512
$ws.Range("I2").Value = 1.289826817062988
$ws.Range("L2").Value = 0.27948156601785
$ws.Range("N2").Value = 1.540167290790905
$ws.Range("B3").Value = 1.941351996942217
$ws.Range("C3").Value = 0.251536935539491
$ws.Range("D3").Value = 0.07080730756334042
$ws.Range("E3").Value = 0.04344196918483956
$ws.Range("F3").Value = 1.780395070784451
$ws.Range("I3").Value = 1.28186931948224
$ws.Range("L3").Value = 0.2688414106238639
$ws.Range("N3").Value = 1.55967806017194
$ws.Range("B4").Value = 1.857655497803705
$ws.Range("C4").Value = 0.2292199559325354
$ws.Range("D4").Value = 0.07135918046992096
$ws.Range("E4").Value = 0.04356372321497792
$ws.Range("F4").Value = 1.759561681787233
$ws.Range("I4").Value = 1.277813728639636
$ws.Range("L4").Value = 0.2624823815995825
$ws.Range("N4").Value = 1.57228175958884
$ws.Range("B5").Value = 1.823828780032329
$ws.Range("C5").Value = 0.2201354367018382
$ws.Range("D5").Value = 0.07159230181207832
$ws.Range("E5").Value = 0.04361546849076436
$ws.Range("F5").Value = 1.751391497944354
$ws.Range("I5").Value = 1.276368864964603
$ws.Range("L5").Value = 0.2599345921770606
$ws.Range("N5").Value = 1.577574367419842
$ws.Range("B6").Value = 1.818228762533977
$ws.Range("C6").Value = 0.2186275314162742
$ws.Range("D6").Value = 0.07163150758817771
$ws.Range("E6").Value = 0.04362418947925795
$ws.Range("F6").Value = 1.750054091567804
$ws.Range("I6").Value = 1.276141470115633
$ws.Range("L6").Value = 0.2595141591016414
$ws.Range("N6").Value = 1.578462640997991
$ws.Range("B7").Value = 1.857198166068656
$ws.Range("C7").Value = 0.2290974000922859
$ws.Range("D7").Value = 0.07136229114764348
$ws.Range("E7").Value = 0.04356441244246767
$ws.Range("F7").Value = 1.759450204382034
$ws.Range("I7").Value = 1.277793402538606
$ws.Range("L7").Value = 0.2624478450875074
$ws.Range("N7").Value = 1.572352504460572
$ws.Range("B8").Value = 2.031626506242048
$ws.Range("C8").Value = 0.2753869384102359
$ws.Range("D8").Value = 0.07024661597103332
$ws.Range("E8").Value = 0.04331931523094168
$ws.Range("F8").Value = 1.80372411454502
$ws.Range("I8").Value = 1.286910138741419
$ws.Range("L8").Value = 0.2757765841642907
$ws.Range("N8").Value = 1.546764843106917
$ws.Range("B9").Value = 2.382666600761866
$ws.Range("C9").Value = 0.3665324717059661
$ws.Range("D9").Value = 0.06832062168480491
$ws.Range("E9").Value = 0.04290560281571087
$ws.Range("F9").Value = 1.900785110235091
$ws.Range("I9").Value = 1.311423332185981
$ws.Range("L9").Value = 0.3033064193138415
$ws.Range("N9").Value = 1.501556601929686
$ws.Range("B10").Value = 2.64621296596988
$ws.Range("C10").Value = 0.4337833333545404
$ws.Range("D10").Value = 0.067066973191654
$ws.Range("E10").Value = 0.04264223187753302
$ws.Range("F10").Value = 1.978480094564105
$ws.Range("I10").Value = 1.333547403678693
$ws.Range("L10").Value = 0.3243992698671576
$ws.Range("N10").Value = 1.471394963443281
$ws.Range("B11").Value = 2.767367854810516
$ws.Range("C11").Value = 0.4644530178136961
$ws.Range("D11").Value = 0.06653208560950929
$ws.Range("E11").Value = 0.04253118571275438
$ws.Range("F11").Value = 2.015243420357422
$ws.Range("I11").Value = 1.34452095919788
$ws.Range("L11").Value = 0.3341872422333694
$ws.Range("N11").Value = 1.458341714217561
$ws.Range("B12").Value = 2.813430589062705
$ws.Range("C12").Value = 0.4760788523372526
$ws.Range("D12").Value = 0.06633466289549617
$ws.Range("E12").Value = 0.04249039213728611
$ws.Range("F12").Value = 2.02937126076128
$ws.Range("I12").Value = 1.348808300675657
$ws.Range("L12").Value = 0.3379216943105092
$ws.Range("N12").Value = 1.453495244174846
$ws.Range("B13").Value = 2.803501940289948
$ws.Range("C13").Value = 0.473574478696321
$ws.Range("D13").Value = 0.06637695276752709
$ws.Range("E13").Value = 0.04249912189430072
$ws.Range("F13").Value = 2.026319363222683
$ws.Range("I13").Value = 1.347879060434551
$ws.Range("L13").Value = 0.3371161657773172
$ws.Range("N13").Value = 1.454534717163703
$ws.Range("B14").Value = 2.771153764760925
$ws.Range("C14").Value = 0.4654092391189693
$ws.Range("D14").Value = 0.06651574054782472
$ws.Range("E14").Value = 0.04252780442077908
$ws.Range("F14").Value = 2.01640157741474
$ws.Range("I14").Value = 1.344871031353065
$ws.Range("L14").Value = 0.3344939160824367
$ws.Range("N14").Value = 1.457941054859177
$ws.Range("B15").Value = 2.75136359403416
$ws.Range("C15").Value = 0.4604093674288947
$ws.Range("D15").Value = 0.06660142093483401
$ws.Range("E15").Value = 0.04254553693720098
$ws.Range("F15").Value = 2.010353589485391
$ws.Range("I15").Value = 1.34304573765057
$ws.Range("L15").Value = 0.332891360931356
$ws.Range("N15").Value = 1.460040119857773
$ws.Range("B16").Value = 2.638320852910397
$ws.Range("C16").Value = 0.4317806158774147
$ws.Range("D16").Value = 0.06710264487506379
$ws.Range("E16").Value = 0.04264966508599599
$ws.Range("F16").Value = 1.97610627978159
$ws.Range("I16").Value = 1.332848649618441
$ws.Range("L16").Value = 0.3237635018914204
$ws.Range("N16").Value = 1.472261500037078
$ws.Range("B17").Value = 2.569298471939192
$ws.Range("C17").Value = 0.4142380929308729
$ws.Range("D17").Value = 0.06741922303091386
$ws.Range("E17").Value = 0.04271578659201403
$ws.Range("F17").Value = 1.955461817716412
$ws.Range("I17").Value = 1.326826749837593
$ws.Range("L17").Value = 0.3182133913064007
$ws.Range("N17").Value = 1.479930244748978
$ws.Range("B18").Value = 2.529717698968682
$ws.Range("C18").Value = 0.4041552884346515
$ws.Range("D18").Value = 0.06760464051449233
$ws.Range("E18").Value = 0.04275464287704489
$ws.Range("F18").Value = 1.943721141895679
$ws.Range("I18").Value = 1.323448621175913
$ws.Range("L18").Value = 0.3150392400649054
$ws.Range("N18").Value = 1.484403893059408
$ws.Range("B19").Value = 2.516336723212532
$ws.Range("C19").Value = 0.4007426340755842
$ws.Range("D19").Value = 0.06766799053280792
$ws.Range("E19").Value = 0.04276794072036738
$ws.Range("F19").Value = 1.93976880691892
$ws.Range("I19").Value = 1.322319496233206
$ws.Range("L19").Value = 0.3139676326933056
$ws.Range("N19").Value = 1.485929361408282
$ws.Range("B20").Value = 2.576633684742148
$ws.Range("C20").Value = 0.4161047742753112
$ws.Range("D20").Value = 0.06738517782971343
$ws.Range("E20").Value = 0.04270866248730387
$ws.Range("F20").Value = 1.957645625857737
$ws.Range("I20").Value = 1.327458933222658
$ws.Range("L20").Value = 0.3188023323659763
$ws.Range("N20").Value = 1.47910739148827
$ws.Range("B21").Value = 2.780650198641695
$ws.Range("C21").Value = 0.4678072385381142
$ws.Range("D21").Value = 0.06647483575439139
$ws.Range("E21").Value = 0.04251934556994108
$ws.Range("F21").Value = 2.019309052722832
$ws.Range("I21").Value = 1.34575097358119
$ws.Range("L21").Value = 0.3352633735926531
$ws.Range("N21").Value = 1.456937907273108
$ws.Range("B22").Value = 2.915060244948336
$ws.Range("C22").Value = 0.5016674892333981
$ws.Range("D22").Value = 0.0659097797028565
$ws.Range("E22").Value = 0.04240294265120736
$ws.Range("F22").Value = 2.060813595457432
$ws.Range("I22").Value = 1.358475205469475
$ws.Range("L22").Value = 0.3461847018054556
$ws.Range("N22").Value = 1.443011669955865
$ws.Range("B23").Value = 2.843224206348452
$ws.Range("C23").Value = 0.483588981447383
$ws.Range("D23").Value = 0.06620861213720985
$ws.Range("E23").Value = 0.04246439961763215
$ws.Range("F23").Value = 2.038550929903266
$ws.Range("I23").Value = 1.351613265308401
$ws.Range("L23").Value = 0.3403407786029931
$ws.Range("N23").Value = 1.450392688805557
$ws.Range("B24").Value = 2.573317120231934
$ws.Range("C24").Value = 0.4152608395239099
$ws.Range("D24").Value = 0.06740055903885533
$ws.Range("E24").Value = 0.04271188067164777
$ws.Range("F24").Value = 1.956657927362244
$ws.Range("I24").Value = 1.327172861741317
$ws.Range("L24").Value = 0.3185360201729708
$ws.Range("N24").Value = 1.479479201697913
$ws.Range("B25").Value = 2.286722651239074
$ws.Range("C25").Value = 0.3418292674201098
$ws.Range("D25").Value = 0.06881344343196361
$ws.Range("E25").Value = 0.04301038084995223
$ws.Range("F25").Value = 1.873416736566327
$ws.Range("I25").Value = 1.304074413376455
$ws.Range("L25").Value = 0.2957079857955875
$ws.Range("N25").Value = 1.51325225136075
